# Generate Report for Handoff
# Adds a new tracked file (8b450931-...) as "Ready for handoff" to all
# three worksheets (Overview, zh-cn, de-de), mirroring the existing row
# for the 1c2b648b-... file.

$wb = $excel.ActiveWorkbook

$mdFileName   = "8b450931-4e5e-4786-86c4-2dc511817c8cooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$e2eMdPath    = "e2e\8b450931-4e5e-4786-86c4-2dc511817c8cooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$readyStatus  = "Ready for handoff"
$dtHandoff    = "2016-08-21 06:36:32"
$dtZhHandoff  = "2016-08-21 06:36:28"
$epochDate    = "0001-01-01 00:00:00"
$zhCnXlf      = "8b450931-4e5e-4786-86c4-2dc511817c8coooooooooooooooooooooooooooooooooooooooo.7039a5c5cddaefc20f2772423356e20ee06ab553.zh-cn.xlf"
$deDeXlf      = "8b450931-4e5e-4786-86c4-2dc511817c8coooooooooooooooooooooooooooooooooooooooo.7039a5c5cddaefc20f2772423356e20ee06ab553.de-de.xlf"
$hyperlinkUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4912e4e6e364018b0e2b702b0c577abdadfd6050/e2e/8b450931-4e5e-4786-86c4-2dc511817c8cooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $mdFileName
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $readyStatus
$wsOverview.Range("F3").Value = $readyStatus
$wsOverview.Range("G3").Value = $dtHandoff

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkUrl, "", "", $e2eMdPath) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $readyStatus
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $zhCnXlf
$wsZhCn.Range("H3").Value = $dtZhHandoff
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = $epochDate
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkUrl, "", "", $mdFileName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $readyStatus
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $deDeXlf
$wsDeDe.Range("H3").Value = $dtHandoff
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = $epochDate
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkUrl, "", "", $mdFileName) | Out-Null

Write-Host "Report generated for handoff."
